$wb = $excel.ActiveWorkbook

function Set-Cell($ws, [string]$addr, $val) {
    $ws.Range($addr).Value = $val
}

$ws = $wb.Worksheets.Item('Citywide Totals')
Set-Cell $ws 'J2' 7643
Set-Cell $ws 'J3' 8024
Set-Cell $ws 'I4' 1777
Set-Cell $ws 'J4' 1746
Set-Cell $ws 'J6' 10972
Set-Cell $ws 'I7' 26232
Set-Cell $ws 'J7' 29008

$ws = $wb.Worksheets.Item('By Neighborhood')
Set-Cell $ws 'J6' 221
Set-Cell $ws 'J7' 828
Set-Cell $ws 'J8' 1837
Set-Cell $ws 'J10' 206
Set-Cell $ws 'J11' 528
Set-Cell $ws 'J15' 357
Set-Cell $ws 'J19' 847
Set-Cell $ws 'J20' 628
Set-Cell $ws 'J23' 266
Set-Cell $ws 'J24' 105
Set-Cell $ws 'J26' 57
Set-Cell $ws 'J29' 1548
Set-Cell $ws 'J30' 101
Set-Cell $ws 'J33' 1305
Set-Cell $ws 'J34' 133
Set-Cell $ws 'J36' 396
Set-Cell $ws 'J37' 894
Set-Cell $ws 'J40' 65
Set-Cell $ws 'J43' 241
Set-Cell $ws 'J44' 231
Set-Cell $ws 'J47' 209
Set-Cell $ws 'J50' 182
Set-Cell $ws 'J52' 738
Set-Cell $ws 'J53' 436
Set-Cell $ws 'J54' 565
Set-Cell $ws 'I63' 186
Set-Cell $ws 'J63' 86
Set-Cell $ws 'J65' 729
Set-Cell $ws 'J67' 1055
Set-Cell $ws 'J75' 85
Set-Cell $ws 'J79' 796
Set-Cell $ws 'J83' 589
Set-Cell $ws 'J84' 239
Set-Cell $ws 'J85' 1187
Set-Cell $ws 'J87' 98
Set-Cell $ws 'J88' 308
Set-Cell $ws 'J89' 364
Set-Cell $ws 'J90' 305
Set-Cell $ws 'J91' 333
Set-Cell $ws 'J92' 96
Set-Cell $ws 'J94' 322
Set-Cell $ws 'J95' 413
Set-Cell $ws 'J96' 326
Set-Cell $ws 'J99' 442
Set-Cell $ws 'I101' 26232
Set-Cell $ws 'J101' 29008

$ws = $wb.Worksheets.Item('West Ridge')
Set-Cell $ws 'J4' 19
Set-Cell $ws 'J7' 326

$ws = $wb.Worksheets.Item('Auburn Gresham')
Set-Cell $ws 'J3' 252
Set-Cell $ws 'J6' 263
Set-Cell $ws 'J7' 828

$ws = $wb.Worksheets.Item('Belmont Cragin')
Set-Cell $ws 'J6' 254
Set-Cell $ws 'J7' 528

$ws = $wb.Worksheets.Item('Uptown')
Set-Cell $ws 'J2' 110
Set-Cell $ws 'J3' 101
Set-Cell $ws 'J6' 114
Set-Cell $ws 'J7' 364

$ws = $wb.Worksheets.Item('South Shore')
Set-Cell $ws 'J2' 315
Set-Cell $ws 'J3' 431
Set-Cell $ws 'J6' 340
Set-Cell $ws 'J7' 1187

$ws = $wb.Worksheets.Item('Little Village')
Set-Cell $ws 'J2' 177
Set-Cell $ws 'J7' 738

$ws = $wb.Worksheets.Item('Logan Square')
Set-Cell $ws 'J6' 284
Set-Cell $ws 'J7' 436

$ws = $wb.Worksheets.Item('Austin')
Set-Cell $ws 'J2' 484
Set-Cell $ws 'J3' 526
Set-Cell $ws 'J4' 95
Set-Cell $ws 'J7' 1837

$ws = $wb.Worksheets.Item('South Chicago')
Set-Cell $ws 'J2' 177
Set-Cell $ws 'J3' 214
Set-Cell $ws 'J7' 589

$ws = $wb.Worksheets.Item('Garfield Park')
Set-Cell $ws 'J6' 467
Set-Cell $ws 'J7' 1305

$ws = $wb.Worksheets.Item('West Pullman')
Set-Cell $ws 'J2' 148
Set-Cell $ws 'J6' 86
Set-Cell $ws 'J7' 413

$ws = $wb.Worksheets.Item('Grand Crossing')
Set-Cell $ws 'J3' 300
Set-Cell $ws 'J4' 33
Set-Cell $ws 'J6' 260
Set-Cell $ws 'J7' 894

$ws = $wb.Worksheets.Item('New City')
Set-Cell $ws 'J2' 212
Set-Cell $ws 'J6' 274
Set-Cell $ws 'J7' 729

$ws = $wb.Worksheets.Item('Woodlawn')
Set-Cell $ws 'J2' 121
Set-Cell $ws 'J6' 117
Set-Cell $ws 'J7' 442

$ws = $wb.Worksheets.Item('Fuller Park')
Set-Cell $ws 'J2' 37
Set-Cell $ws 'J7' 101

$ws = $wb.Worksheets.Item('North Lawndale')
Set-Cell $ws 'J2' 271
Set-Cell $ws 'J3' 397
Set-Cell $ws 'J4' 69
Set-Cell $ws 'J7' 1055

$ws = $wb.Worksheets.Item('South Deering')
Set-Cell $ws 'J3' 76
Set-Cell $ws 'J7' 239

$ws = $wb.Worksheets.Item('Loop')
Set-Cell $ws 'J3' 113
Set-Cell $ws 'J4' 45
Set-Cell $ws 'J6' 259
Set-Cell $ws 'J7' 565

$ws = $wb.Worksheets.Item('Englewood')
Set-Cell $ws 'J3' 546
Set-Cell $ws 'J6' 394
Set-Cell $ws 'J7' 1548

$ws = $wb.Worksheets.Item('Chatham')
Set-Cell $ws 'J2' 211
Set-Cell $ws 'J3' 234
Set-Cell $ws 'J4' 43
Set-Cell $ws 'J7' 847

$ws = $wb.Worksheets.Item('Irving Park')
Set-Cell $ws 'J2' 71
Set-Cell $ws 'J6' 93
Set-Cell $ws 'J7' 231

$ws = $wb.Worksheets.Item('Ashburn')
Set-Cell $ws 'J4' 16
Set-Cell $ws 'J7' 221

$ws = $wb.Worksheets.Item('Avondale')
Set-Cell $ws 'J3' 35
Set-Cell $ws 'J7' 206

$ws = $wb.Worksheets.Item('Dunning')
Set-Cell $ws 'J6' 31
Set-Cell $ws 'J7' 105

$ws = $wb.Worksheets.Item('Douglas')
Set-Cell $ws 'J2' 74
Set-Cell $ws 'J6' 73
Set-Cell $ws 'J7' 266

$ws = $wb.Worksheets.Item('Washington Park')
Set-Cell $ws 'J2' 91
Set-Cell $ws 'J7' 333

$ws = $wb.Worksheets.Item('Roseland')
Set-Cell $ws 'J2' 231
Set-Cell $ws 'J7' 796

$ws = $wb.Worksheets.Item('Chicago Lawn')
Set-Cell $ws 'J2' 172
Set-Cell $ws 'J7' 628

$ws = $wb.Worksheets.Item('Grand Boulevard')
Set-Cell $ws 'J3' 130
Set-Cell $ws 'J6' 118
Set-Cell $ws 'J7' 396

$ws = $wb.Worksheets.Item('Garfield Ridge')
Set-Cell $ws 'J3' 36
Set-Cell $ws 'J7' 133

$ws = $wb.Worksheets.Item('West Loop')
Set-Cell $ws 'J3' 62
Set-Cell $ws 'J6' 172
Set-Cell $ws 'J7' 322

$ws = $wb.Worksheets.Item('Kenwood')
Set-Cell $ws 'J6' 96
Set-Cell $ws 'J7' 209

$ws = $wb.Worksheets.Item('Brighton Park')
Set-Cell $ws 'J6' 166
Set-Cell $ws 'J7' 357

$ws = $wb.Worksheets.Item('Lincoln Square')
Set-Cell $ws 'J2' 48
Set-Cell $ws 'J6' 62
Set-Cell $ws 'J7' 182

$ws = $wb.Worksheets.Item('East Village')
Set-Cell $ws 'J6' 41
Set-Cell $ws 'J7' 57

$ws = $wb.Worksheets.Item('West Elsdon')
Set-Cell $ws 'J2' 29
Set-Cell $ws 'J7' 96

$ws = $wb.Worksheets.Item('United Center')
Set-Cell $ws 'J2' 65
Set-Cell $ws 'J4' 9
Set-Cell $ws 'J6' 162
Set-Cell $ws 'J7' 308

$ws = $wb.Worksheets.Item('Pullman')
Set-Cell $ws 'J3' 28
Set-Cell $ws 'J7' 85

$ws = $wb.Worksheets.Item('Washington Heights')
Set-Cell $ws 'J2' 109
Set-Cell $ws 'J7' 305

$ws = $wb.Worksheets.Item('Hyde Park')
Set-Cell $ws 'J4' 23
Set-Cell $ws 'J6' 144
Set-Cell $ws 'J7' 241

$ws = $wb.Worksheets.Item('Hegewisch')
Set-Cell $ws 'J2' 24
Set-Cell $ws 'J7' 65

$ws = $wb.Worksheets.Item('Ukrainian Village')
Set-Cell $ws 'J3' 14
Set-Cell $ws 'J7' 98
